$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '24.703.15'
Set-TextValue $ws.Range("E2") '  +1.85%  '

Set-TextValue $ws.Range("D3") '1.708.58'
Set-TextValue $ws.Range("E3") '  +2.20%  '

Set-TextValue $ws.Range("D4") '1.011'
Set-TextValue $ws.Range("E4") '  +0.69%  '

Set-TextValue $ws.Range("D5") '310.32'
Set-TextValue $ws.Range("E5") '  +0.73%  '

Set-TextValue $ws.Range("D6") '1.002'
Set-TextValue $ws.Range("E6") '  +0.25%  '

Set-TextValue $ws.Range("D7") '0.3760'
Set-TextValue $ws.Range("E7") '  +0.71%  '

Set-TextValue $ws.Range("D8") '49.67'
Set-TextValue $ws.Range("E8") '  +4.38%  '

Set-TextValue $ws.Range("D9") '0.3472'
Set-TextValue $ws.Range("E9") '  +1.26%  '

Set-TextValue $ws.Range("D10") '1.186'
Set-TextValue $ws.Range("E10") '  +0.53%  '

Set-TextValue $ws.Range("D11") '0.07446'
Set-TextValue $ws.Range("E11") '  +2.28%  '

Set-TextValue $ws.Range("D12") '1.003'
Set-TextValue $ws.Range("E12") '  +0.28%  '

Set-TextValue $ws.Range("D13") '6.276'
Set-TextValue $ws.Range("E13") '  +2.98%  '

Set-TextValue $ws.Range("D14") '20.78'
Set-TextValue $ws.Range("E14") '  +1.60%  '

Set-TextValue $ws.Range("D15") '6.965'
Set-TextValue $ws.Range("E15") '  +3.24%  '

Set-TextValue $ws.Range("E16") '  +2.66%  '

Set-TextValue $ws.Range("D17") '0.00001119'
Set-TextValue $ws.Range("E17") '  +1.08%  '

Set-TextValue $ws.Range("D18") '1.002'
Set-TextValue $ws.Range("E18") '  +0.19%  '

Set-TextValue $ws.Range("D19") '0.06711'
Set-TextValue $ws.Range("E19") '  +0.03%  '

Set-TextValue $ws.Range("D20") '83.91'
Set-TextValue $ws.Range("E20") '  +2.99%  '

Set-TextValue $ws.Range("D21") '17.16'
Set-TextValue $ws.Range("E21") '  +4.70%  '

Set-TextValue $ws.Range("D22") '6.382'
Set-TextValue $ws.Range("E22") '  +4.02%  '

Set-TextValue $ws.Range("D23") '12.99'
Set-TextValue $ws.Range("E23") '  +7.95%  '

Set-TextValue $ws.Range("D24") '24.795.41'
Set-TextValue $ws.Range("E24") '  +2.46%  '

Set-TextValue $ws.Range("D25") '2.442'
Set-TextValue $ws.Range("E25") '  +1.71%  '

Set-TextValue $ws.Range("D26") '2.780'
Set-TextValue $ws.Range("E26") '  +4.84%  '

Set-TextValue $ws.Range("D27") '20.47'
Set-TextValue $ws.Range("E27") '  +5.17%  '

Set-TextValue $ws.Range("D28") '151.02'
Set-TextValue $ws.Range("E28") '  -0.45%  '

Set-TextValue $ws.Range("D29") '131.90'
Set-TextValue $ws.Range("E29") '  +3.70%  '

Set-TextValue $ws.Range("D30") '1.911.44'
Set-TextValue $ws.Range("E30") '  +2.75%  '

Set-TextValue $ws.Range("D31") '1.176'
Set-TextValue $ws.Range("E31") '  +18.95%  '

Set-TextValue $ws.Range("D32") '6.790'
Set-TextValue $ws.Range("E32") '  +6.70%  '

Set-TextValue $ws.Range("D33") '4.228'
Set-TextValue $ws.Range("E33") '  +4.17%  '

Set-TextValue $ws.Range("B34") 'Stellar'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D34") '0.08857'
Set-TextValue $ws.Range("E34") '  +4.77%  '

Set-TextValue $ws.Range("B35") 'Aptos'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D35") '13.73'
Set-TextValue $ws.Range("E35") '  +9.17%  '

Set-TextValue $ws.Range("D36") '1.768'
Set-TextValue $ws.Range("E36") '  +1.23%  '

Set-TextValue $ws.Range("D37") '5.588'
Set-TextValue $ws.Range("E37") '  +4.26%  '

Set-TextValue $ws.Range("D38") '0.06515'
Set-TextValue $ws.Range("E38") '  +1.02%  '

Set-TextValue $ws.Range("D39") '0.02396'
Set-TextValue $ws.Range("E39") '  +2.20%  '

Set-TextValue $ws.Range("D40") '8.946'
Set-TextValue $ws.Range("E40") '  +1.95%  '

Set-TextValue $ws.Range("B41") 'TrustWalletToken'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D41") '1.283'
Set-TextValue $ws.Range("E41") '  +0.03%  '

Set-TextValue $ws.Range("B42") 'Algorand'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D42") '0.2186'
Set-TextValue $ws.Range("E42") '  +3.73%  '

Set-TextValue $ws.Range("D43") '0.6398'
Set-TextValue $ws.Range("E43") '  +4.00%  '

Set-TextValue $ws.Range("D44") '0.9991'
Set-TextValue $ws.Range("E44") '  +0.01%  '

Set-TextValue $ws.Range("D45") '13.96'
Set-TextValue $ws.Range("E45") '  +6.08%  '

Set-TextValue $ws.Range("D46") '0.6105'
Set-TextValue $ws.Range("E46") '  +2.60%  '

Set-TextValue $ws.Range("D47") '3.821'
Set-TextValue $ws.Range("E47") '  +0.45%  '

Set-TextValue $ws.Range("D48") '2.132'
Set-TextValue $ws.Range("E48") '  +5.68%  '

Set-TextValue $ws.Range("D49") '129.31'
Set-TextValue $ws.Range("E49") '  +1.24%  '

Set-TextValue $ws.Range("D50") '0.07262'
Set-TextValue $ws.Range("E50") '  +1.53%  '

Set-TextValue $ws.Range("D51") '79.77'
Set-TextValue $ws.Range("E51") '  +4.14%  '
